# status tracker update (week 6-7)
# Applies updates to the "Table 2" worksheet (sheet2.xml):
#  - fills in Actual Hours (col E) and marks Status (col F) as "Done" for the
#    Week 6 activities (rows 38-44)
#  - fills in the Week 7 activity rows (46-49) that previously existed only as
#    empty styled placeholder rows
#  - appends two brand-new Week 7 activity rows (50-51)
#  - updates the active selection to B49

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1) Week 6 rows (38-44): record Actual Hours in column E and flip Status
#    (column F) from "Planned" to "Done". Re-use the formatting that is
#    already applied to the "Done" status cell on row 33 (fill + border).
# ---------------------------------------------------------------------------
$ws.Range("F33").Copy()
$ws.Range("F38:F44").PasteSpecial(-4122)
$ws.Range("F38:F44").Value2 = "Done"

$ws.Cells.Item(38, 5).Value2 = 0.1
$ws.Cells.Item(39, 5).Value2 = 0.1
$ws.Cells.Item(40, 5).Value2 = 1.5
$ws.Cells.Item(41, 5).Value2 = 2
$ws.Cells.Item(42, 5).Value2 = 1.5
$ws.Cells.Item(43, 5).Value2 = 2.5
$ws.Cells.Item(44, 5).Value2 = 1

# ---------------------------------------------------------------------------
# 2) Row 46 switches from the thin "blank spacer" look (style 24, the same as
#    rows 47-49) to the thicker bordered look used by row 38 (style 23) and
#    gets populated with a new activity. Copy the formatting from row 38
#    first (before row 38's own F-column formatting changed above would be
#    copied along with it) -- row 38's F cell format was already fixed up in
#    step 1 so it is safe to copy the whole row now.
# ---------------------------------------------------------------------------
$ws.Range("A38:G38").Copy()
$ws.Range("A46:G46").PasteSpecial(-4122)
$ws.Rows.Item(46).RowHeight = 8.55

$ws.Cells.Item(46, 1).Value2 = "Client meeting"
$ws.Cells.Item(46, 2).Value2 = "Coordination"
$ws.Cells.Item(46, 3).Value2 = "All"
$ws.Cells.Item(46, 4).Value2 = 0.5
$ws.Cells.Item(46, 6).Value2 = "Planned"
$ws.Cells.Item(46, 7).Value2 = "Give updates about the work, and next thing we need to do"

# Row 46's F cell must stay styled like the rest of the row (style 23, the
# "Planned" look), not the "Done" look that was copied in from row 38's F
# cell - restore it from F46's sibling (row 38's "E" column) which already
# carries style 23.
$ws.Range("E38").Copy()
$ws.Range("F46").PasteSpecial(-4122)
$ws.Cells.Item(46, 6).Value2 = "Planned"

# ---------------------------------------------------------------------------
# 3) Rows 47-49 were empty placeholder rows that already use style 24 for
#    every cell, which matches the target formatting, so only the values
#    need to be filled in.
# ---------------------------------------------------------------------------
$ws.Cells.Item(47, 1).Value2 = "Status tracker"
$ws.Cells.Item(47, 2).Value2 = "Documentation"
$ws.Cells.Item(47, 3).Value2 = "Atidipt"
$ws.Cells.Item(47, 4).Value2 = 0.5
$ws.Cells.Item(47, 6).Value2 = "Planned"
$ws.Cells.Item(47, 7).Value2 = "Track each part of project"

$ws.Cells.Item(48, 1).Value2 = "Team meeting"
$ws.Cells.Item(48, 2).Value2 = "Preparation"
$ws.Cells.Item(48, 3).Value2 = "All"
$ws.Cells.Item(48, 4).Value2 = 1
$ws.Cells.Item(48, 6).Value2 = "Planned"
$ws.Cells.Item(48, 7).Value2 = "Discussion of what all work is done and what all needs to be done"

# Columns C and G of row 49 introduce brand-new shared strings that are
# entered later (see step 4) so that the shared-strings table ends up with
# the same ordering the original author produced.
$ws.Cells.Item(49, 1).Value2 = "Css Part of the pages"
$ws.Cells.Item(49, 2).Value2 = "Preparation"
$ws.Cells.Item(49, 4).Value2 = 2
$ws.Cells.Item(49, 6).Value2 = "Planned"

# ---------------------------------------------------------------------------
# 4) Rows 50-51 are brand new rows appended below the previous bottom of the
#    table. They pick up the worksheet's default column styling, except for
#    the Type (B) and Status (F) cells, which are styled like the rest of
#    the block (style 24, copied from row 49).
#
#    The new, not-yet-seen strings introduced in rows 49-51 must be typed in
#    the same order the original author entered them so that they land at
#    the same indices in the shared-strings table: A49, A50, C50, A51, C49,
#    G49, G50, G51.
# ---------------------------------------------------------------------------
$ws.Cells.Item(50, 1).Value2 = "Attach backend with the frontend"
$ws.Cells.Item(50, 3).Value2 = "Aditya,Aniket"
$ws.Cells.Item(51, 1).Value2 = "Minor changes in the Page"
$ws.Cells.Item(49, 3).Value2 = "Agrim,Shivam"
$ws.Cells.Item(49, 7).Value2 = "Do CSS part of the pages made"
$ws.Cells.Item(50, 7).Value2 = "Connect the backend from local to mongodb atlas"
$ws.Cells.Item(51, 7).Value2 = "Make some changes in purchase and status tracker page"

$ws.Cells.Item(50, 4).Value2 = 2
$ws.Cells.Item(51, 3).Value2 = "Atidipt"
$ws.Cells.Item(51, 4).Value2 = 1

$ws.Range("B49").Copy()
$ws.Range("B50").PasteSpecial(-4122)
$ws.Range("B51").PasteSpecial(-4122)
$ws.Range("F49").Copy()
$ws.Range("F50").PasteSpecial(-4122)
$ws.Range("F51").PasteSpecial(-4122)

$ws.Cells.Item(50, 2).Value2 = "Preparation"
$ws.Cells.Item(50, 6).Value2 = "Planned"
$ws.Cells.Item(51, 2).Value2 = "Preparation"
$ws.Cells.Item(51, 6).Value2 = "Planned"

# ---------------------------------------------------------------------------
# 5) Update the active selection to match the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B49").Select()
